$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Workbook-level: add new "_FilterDatabase" scoped defined name
#    (mirrors the existing _0 / _0_0 variants already in the book)
# ------------------------------------------------------------------
$ws.Names.Add("_xlnm._FilterDatabase_0_0_0", "=Sheet1!`$B`$1:`$I`$67")

# ------------------------------------------------------------------
# 2. Rows 63-66: column E ("hecho") held the placeholder text "*";
#    replace it with the numeric value 1 (row 67 already had 1).
#    This also drops the now-unused "*" shared string on save.
# ------------------------------------------------------------------
$ws.Range("E63").Value = 1
$ws.Range("E64").Value = 1
$ws.Range("E65").Value = 1
$ws.Range("E66").Value = 1

# ------------------------------------------------------------------
# 3. New summary rows 68-70 below the data (data runs rows 2:67)
# ------------------------------------------------------------------
# Row 68: book count, total chapters
$ws.Cells.Item(68, 5).Value = 66
$ws.Cells.Item(68, 8).Formula = "=SUM(H2:H67)"

# Row 69: "total" - count of finished books, chapters finished (+75 manual adj.)
$ws.Cells.Item(69, 4).Value = "total"
$ws.Cells.Item(69, 5).Formula = "=COUNTIF(E2:E67,1)"
$ws.Cells.Item(69, 8).Formula = "=SUMIF(E2:E67,1,H2:H67)+75"
$ws.Cells.Item(69, 8).Font.Color = 0

# Row 70: "proportion" - finished / total ratios
$ws.Cells.Item(70, 4).Value = "proportion"
$ws.Cells.Item(70, 5).Formula = "=E69/E68"
$ws.Cells.Item(70, 8).Formula = "=H69/H68"

# ------------------------------------------------------------------
# 4. View state: scroll/selection moved from G40 to E23
# ------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E23").Select()
